$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value2 = 1232.9412
$ws.Range("J129").Value2 = 1585
$ws.Range("L129").Value2 = 4755
$ws.Range("N129").Value2 = -14755
$ws.Range("H132").Value2 = 1922.4321
$ws.Range("I132").Value2 = 1079.9595
$ws.Range("J132").Value2 = 10828.571
$ws.Range("K132").Value2 = 3239.8785
$ws.Range("L132").Value2 = 32485.713
$ws.Range("M132").Value2 = -709.8784999999998
$ws.Range("N132").Value2 = -37545.713
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 233.66667
$ws.Range("I4").Value2 = 200.5
$ws.Range("J4").Value2 = 300
$ws.Range("K4").Value2 = 200.5
$ws.Range("L4").Value2 = 300
$ws.Range("M4").Value2 = -84.5
$ws.Range("N4").Value2 = -532
$ws.Range("H74").Value2 = 277110.88
$ws.Range("I74").Value2 = 386047.8
$ws.Range("J74").Value2 = 74799.42999999999
$ws.Range("K74").Value2 = 386047.8
$ws.Range("L74").Value2 = 74799.42999999999
$ws.Range("M74").Value2 = -385173.8
$ws.Range("N74").Value2 = -76547.42999999999
$ws.Range("H77").Value2 = 277110.88
$ws.Range("I77").Value2 = 386047.8
$ws.Range("J77").Value2 = 74799.42999999999
$ws.Range("K77").Value2 = 1930239
$ws.Range("L77").Value2 = 373997.15
$ws.Range("M77").Value2 = -1925871
$ws.Range("N77").Value2 = -382733.15
$ws.Range("H97").Value2 = 1022.2963
$ws.Range("I97").Value2 = 1184.1177
$ws.Range("J97").Value2 = 747.2
$ws.Range("K97").Value2 = 1184.1177
$ws.Range("L97").Value2 = 747.2
$ws.Range("M97").Value2 = -688.1177
$ws.Range("N97").Value2 = -1739.2
$ws.Range("H132").Value2 = 18509.984
$ws.Range("I132").Value2 = 23450.787
$ws.Range("J132").Value2 = 3996.375
$ws.Range("K132").Value2 = 70352.361
$ws.Range("L132").Value2 = 11989.125
$ws.Range("M132").Value2 = -67822.361
$ws.Range("N132").Value2 = -17049.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 3806.6296
$ws.Range("I86").Value2 = 5959
$ws.Range("J86").Value2 = 2326.875
$ws.Range("K86").Value2 = 5959
$ws.Range("L86").Value2 = 2326.875
$ws.Range("M86").Value2 = -4836
$ws.Range("N86").Value2 = -4572.875
$ws.Range("H89").Value2 = 3806.6296
$ws.Range("I89").Value2 = 5959
$ws.Range("J89").Value2 = 2326.875
$ws.Range("K89").Value2 = 29795
$ws.Range("L89").Value2 = 11634.375
$ws.Range("M89").Value2 = -24179
$ws.Range("N89").Value2 = -22866.375
$ws.Range("H94").Value2 = 1053.2285
$ws.Range("I94").Value2 = 765.95
$ws.Range("J94").Value2 = 1436.2667
$ws.Range("K94").Value2 = 765.95
$ws.Range("L94").Value2 = 1436.2667
$ws.Range("M94").Value2 = -314.95
$ws.Range("N94").Value2 = -2338.2667
$ws.Range("H134").Value2 = 1711.4314
$ws.Range("I134").Value2 = 1204.3572
$ws.Range("J134").Value2 = 4077.7778
$ws.Range("K134").Value2 = 3613.0716
$ws.Range("L134").Value2 = 12233.3334
$ws.Range("M134").Value2 = -1078.0716
$ws.Range("N134").Value2 = -17303.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 42.3
$ws.Range("I7").Value2 = 35.88889
$ws.Range("J7").Value2 = 100
$ws.Range("K7").Value2 = 35.88889
$ws.Range("L7").Value2 = 100
$ws.Range("M7").Value2 = 77.11111
$ws.Range("N7").Value2 = -326
$ws.Range("H58").Value2 = 4662.4595
$ws.Range("I58").Value2 = 6628.4443
$ws.Range("J58").Value2 = 2799.9473
$ws.Range("K58").Value2 = 6628.4443
$ws.Range("L58").Value2 = 2799.9473
$ws.Range("M58").Value2 = -6425.4443
$ws.Range("N58").Value2 = -3205.9473
$ws.Range("H62").Value2 = 1855190
$ws.Range("I62").Value2 = 5558291.5
$ws.Range("J62").Value2 = 3639.15
$ws.Range("K62").Value2 = 5558291.5
$ws.Range("L62").Value2 = 3639.15
$ws.Range("M62").Value2 = -5557667.5
$ws.Range("N62").Value2 = -4887.15
$ws.Range("H65").Value2 = 1855190
$ws.Range("I65").Value2 = 5558291.5
$ws.Range("J65").Value2 = 3639.15
$ws.Range("K65").Value2 = 27791457.5
$ws.Range("L65").Value2 = 18195.75
$ws.Range("M65").Value2 = -27788337.5
$ws.Range("N65").Value2 = -24435.75
$ws.Range("H132").Value2 = 1930.8392
$ws.Range("I132").Value2 = 971.6667
$ws.Range("K132").Value2 = 2915.0001
$ws.Range("M132").Value2 = -385.0001000000002
$ws.Range("H134").Value2 = 1522.6888
$ws.Range("I134").Value2 = 856.46875
$ws.Range("J134").Value2 = 3162.6155
$ws.Range("K134").Value2 = 2569.40625
$ws.Range("L134").Value2 = 9487.8465
$ws.Range("M134").Value2 = -34.40625
$ws.Range("N134").Value2 = -14557.8465
$ws.Range("H136").Value2 = 4662.4595
$ws.Range("I136").Value2 = 6628.4443
$ws.Range("J136").Value2 = 2799.9473
$ws.Range("K136").Value2 = 19885.3329
$ws.Range("L136").Value2 = 8399.841899999999
$ws.Range("M136").Value2 = -17335.3329
$ws.Range("N136").Value2 = -13499.8419
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 65.066666
$ws.Range("I2").Value2 = 8.444445
$ws.Range("J2").Value2 = 150
$ws.Range("K2").Value2 = 50.66667
$ws.Range("L2").Value2 = 900
$ws.Range("M2").Value2 = 62.33333
$ws.Range("N2").Value2 = -1126
$ws.Range("H9").Value2 = 45000470
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 45000470
$ws.Range("K9").Value2 = 0
$ws.Range("L9").Value2 = 135001410
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value2 = -135001858
$ws.Range("H16").Value2 = 566.6667
$ws.Range("I16").Value2 = 566.6667
$ws.Range("K16").Value2 = 1700.0001
$ws.Range("M16").Value2 = -1527.0001
$ws.Range("H17").Value2 = 947.8333
$ws.Range("I17").Value2 = 532
$ws.Range("J17").Value2 = 1363.6666
$ws.Range("K17").Value2 = 1596
$ws.Range("L17").Value2 = 4090.9998
$ws.Range("M17").Value2 = -1427
$ws.Range("N17").Value2 = -4428.9998
$ws.Range("H34").Value2 = 589.4483
$ws.Range("I34").Value2 = 52
$ws.Range("J34").Value2 = 598.8772
$ws.Range("K34").Value2 = 156
$ws.Range("L34").Value2 = 1796.6316
$ws.Range("M34").Value2 = -72
$ws.Range("N34").Value2 = -1964.6316
$ws.Range("H39").Value2 = 3200.3
$ws.Range("J39").Value2 = 3812.5
$ws.Range("L39").Value2 = 11437.5
$ws.Range("N39").Value2 = -12025.5
$ws.Range("H55").Value2 = 2965.5264
$ws.Range("I55").Value2 = 450
$ws.Range("J55").Value2 = 3105.2778
$ws.Range("K55").Value2 = 1350
$ws.Range("L55").Value2 = 9315.8334
$ws.Range("M55").Value2 = -1173
$ws.Range("N55").Value2 = -9669.8334
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 6056.8184
$ws.Range("I70").Value2 = 4492.857
$ws.Range("J70").Value2 = 6786.6665
$ws.Range("K70").Value2 = 4492.857
$ws.Range("L70").Value2 = 6786.6665
$ws.Range("M70").Value2 = -4222.857
$ws.Range("N70").Value2 = -7326.6665
$ws.Range("H73").Value2 = 6056.8184
$ws.Range("I73").Value2 = 4492.857
$ws.Range("J73").Value2 = 6786.6665
$ws.Range("K73").Value2 = 4492.857
$ws.Range("L73").Value2 = 6786.6665
$ws.Range("M73").Value2 = -3556.857
$ws.Range("N73").Value2 = -8658.666499999999
$ws.Range("H132").Value2 = 2928.8103
$ws.Range("I132").Value2 = 2546.0977
$ws.Range("K132").Value2 = 7638.293099999999
$ws.Range("M132").Value2 = -5108.293099999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value2 = 43481644
$ws.Range("I100").Value2 = 3941.1765
$ws.Range("J100").Value2 = 166668460
$ws.Range("K100").Value2 = 3941.1765
$ws.Range("L100").Value2 = 166668460
$ws.Range("M100").Value2 = -3400.1765
$ws.Range("N100").Value2 = -166669542
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 1418.6
$ws.Range("I100").Value2 = 1548.25
$ws.Range("J100").Value2 = 900
$ws.Range("K100").Value2 = 3096.5
$ws.Range("L100").Value2 = 1800
$ws.Range("M100").Value2 = -2555.5
$ws.Range("N100").Value2 = -2882
$ws.Range("H132").Value2 = 1464.6731
$ws.Range("I132").Value2 = 933.72095
$ws.Range("J132").Value2 = 4001.4443
$ws.Range("K132").Value2 = 2801.16285
$ws.Range("L132").Value2 = 12004.3329
$ws.Range("M132").Value2 = -271.1628500000002
$ws.Range("N132").Value2 = -17064.3329
